$d = $word.ActiveDocument

# 1) In the charges table, the "Plea" row has "No Contest" in both data
#    columns. The second column (TAIL LIGHTS-REAR LICENSE PLATE charge)
#    plea changes to "Guilty". Target the specific table cell directly
#    (rather than a document-wide Find/Replace) since the text
#    "No Contest" is not unique in the document.
$tbl = $d.Tables.Item(1)
$pleaCell = $tbl.Cell(4, 3)
$cellStart = $pleaCell.Range.Start
$pleaRange = $d.Range($cellStart, $cellStart + 10)
if ($pleaRange.Text -eq "No Contest") {
    $pleaRange.Text = "Guilty"
}

# 2) License Suspension paragraph: "hunting" -> "driving"
$d.Content.Find.Execute("hunting", $true, $false, $false, $false, $false, `
    $true, 1, $false, "driving", 2) | Out-Null

# 3) License Suspension paragraph: term changes from 12 months to 6 months
$d.Content.Find.Execute("for a term of 12 months", $true, $false, $false, `
    $false, $false, $true, 1, $false, "for a term of 6 months", 2) | Out-Null
